$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.031.26'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.29%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.816.94'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.40%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '337.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.48%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4269'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +11.77%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3511'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.73%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '45.52'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.22%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.147'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.15%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07444'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.99%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.95'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.25%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.002'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.01%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.261'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.88%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.310'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.96%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.814.07'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.07%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001085'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.89%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06699'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.68%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '82.13'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.44%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.08%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.24'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.66%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.409'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.24%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.034.75'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.36%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.85'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.66%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.395'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.29%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.470'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.87%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.73'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.17%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '155.25'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.73%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.019.72'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.06%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.303'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -9.57%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '132.50'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.50%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.075'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.85%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.961'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.25%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.09242'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.57%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.39'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.25%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02376'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.09%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6704'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.78%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06281'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.17%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.230'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.47%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2173'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.75%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.499'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.34%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.221'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.35%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.106'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.00%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '14.25'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.86%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.874'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.29%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.6138'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.93%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '128.11'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.43%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.050'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.81%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.181'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.85%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07112'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.31%  '
